# Bolzenlagerung.pptx - "Finale-Ordner erstellt, kleinere Änderungen"
#
# 1) Update the cached "datetimeFigureOut" date field text on the slide
#    master and on every slide layout (06.11.2019 -> 20.11.2019).
# 2) Update the "140" -> "315" value textbox on slide 1.

$p = $ppt.ActivePresentation

function Update-DateField {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Datumsplatzhalter*" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "06.11.2019") {
                $tr.Text = "20.11.2019"
            }
        }
    }
}

# Slide master date placeholder.
Update-DateField $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateField $layout.Shapes
}

# Slide 1: "140" -> "315" in the value textbox ("Textfeld 37").
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "140") {
        $sh.TextFrame.TextRange.Text = "315"
    }
}
